$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B (rows 2-97), MW production values
$newB = @(1413,1338,1253,1222,1194,1220,1251,1259,1291,1326,1367,1405,1434,1425,1447,1485,1535,1443,1419,1378,1333,1301,1306,1348,1317,1272,1275,1266,1223,1163,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

# Shift all timestamps in column A (rows 2-97) forward by 22 days
# and set the corresponding new production values in column B
for ($i = 0; $i -lt $newB.Length; $i++) {
    $row = $i + 2
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value = $cellA.Value2 + 22
    $ws.Cells.Item($row, 2).Value = $newB[$i]
}
